$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text tweaks (shared strings used by A8 "Volume 30   Number  42"
#    and C9 "Report Covering the Week  10/16/2023  Through  10/22/2023")
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/23/2023  Through  10/29/2023"

# ---------------------------------------------------------------------------
# 2) Cells that flip between a number and the literal text placeholders
#    ("0" = shared string 20, "***.*" = shared string 21). Grab the
#    formatting from a same-styled neighbour cell (column A on the same
#    row carries style 14) via PasteSpecial-formats, then pull the actual
#    text value from a cell that already holds the right shared string so
#    the written cell becomes a genuine t="s" reference instead of a
#    numeric literal.
# ---------------------------------------------------------------------------
function Set-TextPlaceholder($targetAddr, $styleSourceAddr, $valueSourceAddr) {
    $ws.Range($styleSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($valueSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# C16: 1 -> "0"
Set-TextPlaceholder "C16" "A16" "C22"
# D20: 4 -> "0"
Set-TextPlaceholder "D20" "A20" "D14"
# E20: -25 -> "***.*"
Set-TextPlaceholder "E20" "A20" "E14"

$ws.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Plain numeric updates for rows 16-21, 24, 25, 27
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -66.666666666666
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = -27.142857142857
$ws.Range("L16").Value = 155
$ws.Range("M16").Value = -21.538461538461
$ws.Range("N16").Value = -80.970149253731

$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 80
$ws.Range("J17").Value = 79
$ws.Range("K17").Value = 1.265822784810
$ws.Range("L17").Value = 33.333333333333
$ws.Range("M17").Value = 73.913043478260
$ws.Range("N17").Value = -18.367346938775

$ws.Range("C18").Value = 9
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = 8
$ws.Range("I18").Value = 256
$ws.Range("J18").Value = 222
$ws.Range("K18").Value = 15.315315315315
$ws.Range("L18").Value = 43.016759776536
$ws.Range("M18").Value = 25.490196078431
$ws.Range("N18").Value = -69.846878680800

$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = -27.419354838709
$ws.Range("I19").Value = 524
$ws.Range("J19").Value = 542
$ws.Range("K19").Value = -3.321033210332
$ws.Range("L19").Value = 71.241830065359
$ws.Range("M19").Value = 73.509933774834
$ws.Range("N19").Value = 11.965811965812

$ws.Range("C20").Value = 4
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 212.5
$ws.Range("I20").Value = 160
$ws.Range("K20").Value = 72.043010752688
$ws.Range("L20").Value = 180.701754385965
$ws.Range("M20").Value = 35.593220338983
$ws.Range("N20").Value = -94.130594277329

$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -3.333333333333
$ws.Range("F21").Value = 108
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -1.818181818181
$ws.Range("I21").Value = 1081
$ws.Range("J21").Value = 1008
$ws.Range("K21").Value = 7.242063492063
$ws.Range("L21").Value = 72.408293460925
$ws.Range("M21").Value = 46.081081081081
$ws.Range("N21").Value = -75.559574949129

$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 64
$ws.Range("H24").Value = -15.625
$ws.Range("I24").Value = 487
$ws.Range("J24").Value = 649
$ws.Range("K24").Value = -24.961479198767
$ws.Range("L24").Value = 8.705357142857
$ws.Range("M24").Value = 30.913978494623

$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = -33.333333333333
$ws.Range("I25").Value = 191
$ws.Range("J25").Value = 206
$ws.Range("K25").Value = -7.281553398058
$ws.Range("L25").Value = 55.284552845528
$ws.Range("M25").Value = 31.724137931034

$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
